$wb = $excel.ActiveWorkbook

$wsActivities = $wb.Worksheets.Item("Activities")

# Add the "Invoices" worksheet right after "Activities" (the current last sheet)
$wsInvoices = $wb.Worksheets.Add($null, $wsActivities)
$wsInvoices.Name = "Invoices"

$wsInvoices.Range("A1").Value = "Invoice Number"
$wsInvoices.Range("A2").Value = "Inv 1"
$wsInvoices.Range("A3").Value = "Inv 2"
$wsInvoices.Range("B1").Value = "Invoice"
$wsInvoices.Range("B2").Value = "Invoice 1"
$wsInvoices.Range("B3").Value = "Invoice 2"
$wsInvoices.Columns.Item(1).ColumnWidth = 16.5

# Add the "InvoiceContacts" worksheet right after "Invoices"
$wsInvoiceContacts = $wb.Worksheets.Add($null, $wsInvoices)
$wsInvoiceContacts.Name = "InvoiceContacts"

$wsInvoiceContacts.Range("A1").Value = "InvoiceNo"
$wsInvoiceContacts.Range("A2").Value = "Inv 2"
$wsInvoiceContacts.Range("A3").Value = "Inv 1"
$wsInvoiceContacts.Range("B1").Value = "ContactID"
$wsInvoiceContacts.Range("B2").Value = 2
$wsInvoiceContacts.Range("B3").Value = 3
$wsInvoiceContacts.Columns.Item(1).ColumnWidth = 13.0

$wsInvoiceContacts.Range("B2").Select()

# Make InvoiceContacts the active sheet/tab
$wsInvoiceContacts.Activate()
